# Journal de travail: add a new row (line 43) to the "Tableau1" table,
# logging the work session for "Affichage des scores en fonction de la
# maquette" (commit: "Réalisation de l'écran du score en fonction de la
# maquette").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The existing table (Tableau1) covers E5:M42. Copy the formatting of the
# last data row (42) down into the new row (43) so the new row inherits the
# same number formats / alignment / wrap styles as the rest of the table.
$ws.Range("E42:M42").Copy()
$ws.Range("E43:M43").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Match the row height used for the other wrapped-text rows in the sheet.
$ws.Range("E43:M43").RowHeight = 43.2

# Date (2021-03-24), start/end time, same Durée formula pattern as the rows
# above it.
$ws.Range("E43").Value = 44279
$ws.Range("F43").Value = 0.57291666666666663
$ws.Range("G43").Value = 0.60416666666666663
$ws.Range("H43").Formula = "=IF(ISBLANK(Tableau1[[#This Row],[Heure Début]]),"""",Tableau1[[#This Row],[Heure fin]]-Tableau1[[#This Row],[Heure Début]])"

$ws.Range("I43").Value = "Développement"
# Set the source URL before the task/description text so new shared
# strings land in the same order as the target workbook.
$ws.Range("M43").Value = "https://www.ltam.lu/cours-c/prg-c42.htm"
$ws.Range("J43").Value = "Affichage des scores"
$ws.Range("K43").Value = "CPNV"
$ws.Range("L43").Value = "Réalisation de l'écran des scores en fonction de la maquette"

# Grow the ListObject (Excel Table) so E5:M43 (incl. the new row + its
# AutoFilter) is covered.
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("E5:M43"))

# Reflect the view state after entering the new row: scrolled down a bit
# and the cursor sitting one row below/at the new entry.
$ws.Activate() | Out-Null
$ws.Range("B38").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 38
$win.ScrollColumn = 2
$ws.Range("L44").Select() | Out-Null
